$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4352.7144
$ws.Range("J19").Value = 6798.923
$ws.Range("L19").Value = 6798.923
$ws.Range("N19").Value = -7148.923
$ws.Range("H62").Value = 4475
$ws.Range("I62").Value = 3966.6667
$ws.Range("K62").Value = 3966.6667
$ws.Range("M62").Value = -3342.6667
$ws.Range("H65").Value = 4475
$ws.Range("I65").Value = 3966.6667
$ws.Range("K65").Value = 19833.3335
$ws.Range("M65").Value = -16713.3335
$ws.Range("H74").Value = 11832.667
$ws.Range("I74").Value = 7398.4
$ws.Range("K74").Value = 7398.4
$ws.Range("M74").Value = -6462.4
$ws.Range("H76").Value = 9119.272000000001
$ws.Range("I76").Value = 8617.571
$ws.Range("K76").Value = 8617.571
$ws.Range("M76").Value = -8302.571
$ws.Range("H77").Value = 11832.667
$ws.Range("I77").Value = 7398.4
$ws.Range("K77").Value = 36992
$ws.Range("M77").Value = -32312
$ws.Range("H79").Value = 9119.272000000001
$ws.Range("I79").Value = 8617.571
$ws.Range("K79").Value = 8617.571
$ws.Range("M79").Value = -7525.571
$ws.Range("H99").Value = 300.23077
$ws.Range("I99").Value = 294.2
$ws.Range("J99").Value = 320.33334
$ws.Range("K99").Value = 882.5999999999999
$ws.Range("L99").Value = 961.0000200000001
$ws.Range("M99").Value = 615.4000000000001
$ws.Range("N99").Value = -3957.00002
$ws.Range("H112").Value = 1529.6471
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1529.6471
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4588.9413
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6804.9413
$ws.Range("H132").Value = 1653.1364
$ws.Range("I132").Value = 1568.4762
$ws.Range("J132").Value = 3431
$ws.Range("K132").Value = 4705.4286
$ws.Range("L132").Value = 10293
$ws.Range("M132").Value = -2175.4286
$ws.Range("N132").Value = -15353
$ws.Range("H135").Value = 1102.5
$ws.Range("I135").Value = 669.8611
$ws.Range("K135").Value = 6028.7499
$ws.Range("M135").Value = -3493.7499
$ws.Range("H137").Value = 1833.0646
$ws.Range("I137").Value = 1603.9166
$ws.Range("J137").Value = 2618.7144
$ws.Range("K137").Value = 4811.7498
$ws.Range("L137").Value = 7856.1432
$ws.Range("M137").Value = -2261.7498
$ws.Range("N137").Value = -12956.1432
$ws.Range("H141").Value = 967.3333
$ws.Range("I141").Value = 1077.4117
$ws.Range("J141").Value = 499.5
$ws.Range("K141").Value = 3232.2351
$ws.Range("L141").Value = 1498.5
$ws.Range("M141").Value = 1947.7649
$ws.Range("N141").Value = -11858.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6111.347
$ws.Range("I32").Value = 4091.7112
$ws.Range("J32").Value = 42464.8
$ws.Range("K32").Value = 4091.7112
$ws.Range("L32").Value = 42464.8
$ws.Range("M32").Value = -3804.7112
$ws.Range("N32").Value = -43038.8
$ws.Range("H74").Value = 63546.965
$ws.Range("I74").Value = 38685.375
$ws.Range("J74").Value = 182882.6
$ws.Range("K74").Value = 38685.375
$ws.Range("L74").Value = 182882.6
$ws.Range("M74").Value = -37811.375
$ws.Range("N74").Value = -184630.6
$ws.Range("H77").Value = 63546.965
$ws.Range("I77").Value = 38685.375
$ws.Range("J77").Value = 182882.6
$ws.Range("K77").Value = 193426.875
$ws.Range("L77").Value = 914413
$ws.Range("M77").Value = -189058.875
$ws.Range("N77").Value = -923149
$ws.Range("H124").Value = 42642.832
$ws.Range("J124").Value = 42642.832
$ws.Range("L124").Value = 42642.832
$ws.Range("N124").Value = -52462.832
$ws.Range("H125").Value = 137607.5
$ws.Range("J125").Value = 137607.5
$ws.Range("L125").Value = 137607.5
$ws.Range("N125").Value = -147447.5
$ws.Range("H132").Value = 8471.440000000001
$ws.Range("I132").Value = 9125.857
$ws.Range("J132").Value = 5035.75
$ws.Range("K132").Value = 27377.571
$ws.Range("L132").Value = 15107.25
$ws.Range("M132").Value = -24847.571
$ws.Range("N132").Value = -20167.25

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 49749.75
$ws.Range("J109").Value = 49749.75
$ws.Range("L109").Value = 49749.75
$ws.Range("N109").Value = -52523.75
$ws.Range("H134").Value = 3388.8
$ws.Range("I134").Value = 2196.52
$ws.Range("K134").Value = 6589.559999999999
$ws.Range("M134").Value = -4054.559999999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2438.2222
$ws.Range("I31").Value = 1620.9
$ws.Range("J31").Value = 2818.372
$ws.Range("K31").Value = 1620.9
$ws.Range("L31").Value = 2818.372
$ws.Range("M31").Value = -1325.9
$ws.Range("N31").Value = -3408.372
$ws.Range("H34").Value = 2438.2222
$ws.Range("I34").Value = 1620.9
$ws.Range("J34").Value = 2818.372
$ws.Range("K34").Value = 1620.9
$ws.Range("L34").Value = 2818.372
$ws.Range("M34").Value = -1418.9
$ws.Range("N34").Value = -3222.372
$ws.Range("H69").Value = 42895.6
$ws.Range("I69").Value = 31493.666
$ws.Range("K69").Value = 31493.666
$ws.Range("M69").Value = -30744.666
$ws.Range("H72").Value = 42895.6
$ws.Range("I72").Value = 31493.666
$ws.Range("K72").Value = 94480.99800000001
$ws.Range("M72").Value = -90736.99800000001
$ws.Range("H93").Value = 25798.924
$ws.Range("I93").Value = 10658.571
$ws.Range("J93").Value = 43462.668
$ws.Range("K93").Value = 10658.571
$ws.Range("L93").Value = 43462.668
$ws.Range("M93").Value = -8786.571
$ws.Range("N93").Value = -47206.668
$ws.Range("H103").Value = 4178.1665
$ws.Range("I103").Value = 4178.1665
$ws.Range("K103").Value = 4178.1665
$ws.Range("M103").Value = -3006.1665

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 248.91667
$ws.Range("I6").Value = 82
$ws.Range("J6").Value = 749.6667
$ws.Range("K6").Value = 246
$ws.Range("L6").Value = 2249.0001
$ws.Range("M6").Value = -133
$ws.Range("N6").Value = -2475.0001
$ws.Range("H63").Value = 24996.666
$ws.Range("J63").Value = 24996.666
$ws.Range("L63").Value = 74989.99800000001
$ws.Range("N63").Value = -76487.99800000001
$ws.Range("H66").Value = 24996.666
$ws.Range("J66").Value = 24996.666
$ws.Range("L66").Value = 224969.994
$ws.Range("N66").Value = -232457.994
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H75").Value = 20201.8
$ws.Range("I75").Value = 1506
$ws.Range("J75").Value = 32665.666
$ws.Range("K75").Value = 4518
$ws.Range("L75").Value = 97996.99800000001
$ws.Range("M75").Value = -3520
$ws.Range("N75").Value = -99992.99800000001
$ws.Range("H78").Value = 20201.8
$ws.Range("I78").Value = 1506
$ws.Range("J78").Value = 32665.666
$ws.Range("K78").Value = 13554
$ws.Range("L78").Value = 293990.994
$ws.Range("M78").Value = -8562
$ws.Range("N78").Value = -303974.994
$ws.Range("H87").Value = 17750
$ws.Range("I87").Value = 17750
$ws.Range("K87").Value = 53250
$ws.Range("M87").Value = -52002
$ws.Range("H90").Value = 17750
$ws.Range("I90").Value = 17750
$ws.Range("K90").Value = 159750
$ws.Range("M90").Value = -153510
$ws.Range("H114").Value = 3884.6667
$ws.Range("J114").Value = 3702
$ws.Range("L114").Value = 11106
$ws.Range("N114").Value = -17614
$ws.Range("H117").Value = 907.1667
$ws.Range("I117").Value = 600
$ws.Range("J117").Value = 968.6
$ws.Range("K117").Value = 1800
$ws.Range("L117").Value = 2905.8
$ws.Range("M117").Value = 1642
$ws.Range("N117").Value = -9789.799999999999
$ws.Range("H121").Value = 14287565
$ws.Range("I121").Value = 129.66667
$ws.Range("J121").Value = 18184138
$ws.Range("K121").Value = 389.00001
$ws.Range("L121").Value = 54552414
$ws.Range("M121").Value = 920.99999
$ws.Range("N121").Value = -54555034
$ws.Range("H132").Value = 3975.5
$ws.Range("I132").Value = 1826
$ws.Range("J132").Value = 6125
$ws.Range("K132").Value = 16434
$ws.Range("L132").Value = 55125
$ws.Range("M132").Value = -13904
$ws.Range("N132").Value = -60185
$ws.Range("H139").Value = 2936.0952
$ws.Range("I139").Value = 2231
$ws.Range("J139").Value = 7166.6665
$ws.Range("K139").Value = 6693
$ws.Range("L139").Value = 21499.9995
$ws.Range("M139").Value = -1553
$ws.Range("N139").Value = -31779.9995
$ws.Range("H140").Value = 1787.4546
$ws.Range("I140").Value = 1729.7142
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 5189.142599999999
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -9.14259999999922
$ws.Range("N140").Value = -19360

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3406.75
$ws.Range("I68").Value = 1957
$ws.Range("J68").Value = 4442.2856
$ws.Range("K68").Value = 1957
$ws.Range("L68").Value = 4442.2856
$ws.Range("M68").Value = -1208
$ws.Range("N68").Value = -5940.2856
$ws.Range("H71").Value = 3406.75
$ws.Range("I71").Value = 1957
$ws.Range("J71").Value = 4442.2856
$ws.Range("K71").Value = 9785
$ws.Range("L71").Value = 22211.428
$ws.Range("M71").Value = -6041
$ws.Range("N71").Value = -29699.428
$ws.Range("H136").Value = 59452.688
$ws.Range("I136").Value = 3416.0833
$ws.Range("J136").Value = 227562.5
$ws.Range("K136").Value = 10248.2499
$ws.Range("L136").Value = 682687.5
$ws.Range("M136").Value = -7698.249899999999
$ws.Range("N136").Value = -687787.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 63792.5
$ws.Range("I88").Value = 15585
$ws.Range("J88").Value = 112000
$ws.Range("K88").Value = 15585
$ws.Range("L88").Value = 112000
$ws.Range("M88").Value = -15179
$ws.Range("N88").Value = -112812
$ws.Range("H91").Value = 63792.5
$ws.Range("I91").Value = 15585
$ws.Range("J91").Value = 112000
$ws.Range("K91").Value = 15585
$ws.Range("L91").Value = 112000
$ws.Range("M91").Value = -14181
$ws.Range("N91").Value = -114808
$ws.Range("H136").Value = 7631.9565
$ws.Range("I136").Value = 9408.071
$ws.Range("J136").Value = 4869.1113
$ws.Range("K136").Value = 28224.213
$ws.Range("L136").Value = 14607.3339
$ws.Range("M136").Value = -25674.213
$ws.Range("N136").Value = -19707.3339
